# Apply updated bank reconciliation data to the "Transaksi" and "Summary" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaksi")
$summary = $wb.Worksheets.Item("Summary")

# Make sure date columns (A and E) on any newly added rows use the same
# date format as the existing rows (style index backing numFmt YYYY-MM-DD).
$ws.Range("A8:A11").NumberFormat = "YYYY-MM-DD"
$ws.Range("E8:E11").NumberFormat = "YYYY-MM-DD"

# --- Row 2 ---
$ws.Range("A2").Value = 45443
$ws.Range("B2").Value = "0,00"
$ws.Range("C2").Value = "0,00"
$ws.Range("D2").Value = "76.939.992,80"
$ws.Range("E2").Value = 45443
$ws.Range("F2").Value = "0,00"
$ws.Range("G2").Value = "0,00"
$ws.Range("H2").Value = "76.939.992,80"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "Opening Balance"

# --- Row 3 ---
$ws.Range("A3").Value = 45447
$ws.Range("B3").Value = "0,00"
$ws.Range("C3").Value = "33.224.480,00"
$ws.Range("D3").Value = "43.715.512,80"
$ws.Range("E3").Value = 45447
$ws.Range("F3").Value = "33.224.480,00"
$ws.Range("G3").Value = "0,00"
$ws.Range("H3").Value = "43.715.512,80"
$ws.Range("I3").Value = "0,00"
$ws.Range("J3").Value = "0,00"
$ws.Range("K3").Value = "Matched"
$ws.Range("L3").Value = "-"

# --- Row 4 ---
$ws.Range("A4").Value = 45453
$ws.Range("B4").Value = "8.250.000,00"
$ws.Range("C4").Value = "0,00"
$ws.Range("D4").Value = "51.965.512,80"
$ws.Range("E4").Value = 45453
$ws.Range("F4").Value = "0,00"
$ws.Range("G4").Value = "8.250.000,00"
$ws.Range("H4").Value = "51.965.512,80"
$ws.Range("I4").Value = "0,00"
$ws.Range("J4").Value = "0,00"
$ws.Range("K4").Value = "Matched"
$ws.Range("L4").Value = "-"

# --- Row 5 ---
$ws.Range("A5").Value = 45453
$ws.Range("B5").Value = "2.233.833,00"
$ws.Range("C5").Value = "0,00"
$ws.Range("D5").Value = "54.199.345,80"
$ws.Range("E5").Value = 45453
$ws.Range("F5").Value = "0,00"
$ws.Range("G5").Value = "2.233.833,00"
$ws.Range("H5").Value = "54.199.345,80"
$ws.Range("I5").Value = "0,00"
$ws.Range("J5").Value = "0,00"
$ws.Range("K5").Value = "Matched"
$ws.Range("L5").Value = "-"

# --- Row 6 ---
$ws.Range("A6").Value = 45461
$ws.Range("B6").Value = "0,00"
$ws.Range("C6").Value = "130.900,00"
$ws.Range("D6").Value = "54.068.445,80"
$ws.Range("E6").Value = 45461
$ws.Range("F6").Value = "130.900,00"
$ws.Range("G6").Value = "0,00"
$ws.Range("H6").Value = "54.068.445,80"
$ws.Range("I6").Value = "0,00"
$ws.Range("J6").Value = "0,00"
$ws.Range("K6").Value = "Matched"
$ws.Range("L6").Value = "-"

# --- Row 7 ---
$ws.Range("A7").Value = 45462
$ws.Range("B7").Value = "0,00"
$ws.Range("C7").Value = "8.025.000,00"
$ws.Range("D7").Value = "46.043.445,80"
$ws.Range("E7").Value = 45462
$ws.Range("F7").Value = "8.025.000,00"
$ws.Range("G7").Value = "0,00"
$ws.Range("H7").Value = "46.043.445,80"
$ws.Range("I7").Value = "0,00"
$ws.Range("J7").Value = "0,00"
$ws.Range("K7").Value = "Matched"
$ws.Range("L7").Value = "-"

# --- Row 8 (new) ---
$ws.Range("A8").Value = 45463
$ws.Range("B8").Value = "0,00"
$ws.Range("C8").Value = "1.881.550,00"
$ws.Range("D8").Value = "44.161.895,80"
$ws.Range("E8").Value = 45463
$ws.Range("F8").Value = "1.881.550,00"
$ws.Range("G8").Value = "0,00"
$ws.Range("H8").Value = "44.161.895,80"
$ws.Range("I8").Value = "0,00"
$ws.Range("J8").Value = "0,00"
$ws.Range("K8").Value = "Matched"
$ws.Range("L8").Value = "-"

# --- Row 9 (new) ---
$ws.Range("A9").Value = 45473
$ws.Range("B9").Value = "10.350,00"
$ws.Range("C9").Value = "0,00"
$ws.Range("D9").Value = "44.172.245,80"
$ws.Range("E9").Value = 45473
$ws.Range("F9").Value = "0,00"
$ws.Range("G9").Value = "10.350,00"
$ws.Range("H9").Value = "44.172.245,80"
$ws.Range("I9").Value = "0,00"
$ws.Range("J9").Value = "0,00"
$ws.Range("K9").Value = "Matched"
$ws.Range("L9").Value = "-"

# --- Row 10 (new) ---
$ws.Range("A10").Value = 45473
$ws.Range("B10").Value = "0,00"
$ws.Range("C10").Value = "2.070,00"
$ws.Range("D10").Value = "44.170.175,80"
$ws.Range("E10").Value = 45473
$ws.Range("F10").Value = "2.070,00"
$ws.Range("G10").Value = "0,00"
$ws.Range("H10").Value = "44.170.175,80"
$ws.Range("I10").Value = "0,00"
$ws.Range("J10").Value = "0,00"
$ws.Range("K10").Value = "Matched"
$ws.Range("L10").Value = "-"

# --- Row 11 (new, Closing Balance) ---
$ws.Range("A11").Value = 45473
$ws.Range("B11").Value = "10.494.183,00"
$ws.Range("C11").Value = "43.264.000,00"
$ws.Range("D11").Value = "44.170.175,80"
$ws.Range("E11").Value = 45473
$ws.Range("F11").Value = "43.264.000,00"
$ws.Range("G11").Value = "10.494.183,00"
$ws.Range("H11").Value = "44.170.175,80"
$ws.Range("I11").Value = "0,00"
$ws.Range("J11").Value = "0,00"
$ws.Range("K11").Value = "Closing Balance"

# --- Summary sheet ---
$summary.Range("B2").Value = 45443
$summary.Range("C2").Value = "76.939.992,80"
$summary.Range("D2").Value = "76.939.992,80"

$summary.Range("B3").Value = 45473
$summary.Range("C3").Value = "44.170.175,80"
$summary.Range("D3").Value = "44.170.175,80"
